# Swap the two occurrence rows (row 2 <-> row 3) and refresh a few
# fields (coordinates get rounded, the Starttid/Sluttid columns are
# dropped) to match the freshly re-exported source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (becomes the former row-3 record) ---
$ws.Range("A2").Value = 112188960
$ws.Range("B2").Value = 89405
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 1202
$ws.Range("F2").Value = "Ullticka"
$ws.Range("G2").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H2").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q2").Value = 332391
$ws.Range("R2").Value = 6627094
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# --- Row 3 (becomes the former row-2 record) ---
$ws.Range("A3").Value = 112188940
$ws.Range("B3").Value = 101703
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 222412
$ws.Range("F3").Value = "Tibast"
$ws.Range("G3").Value = "Daphne mezereum"
$ws.Range("H3").Value = "L."
$ws.Range("Q3").Value = 332308
$ws.Range("R3").Value = 6627086
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
